# The document has two BTec_Logo-Orange images (one in the "default"
# header, one in the "first page" header) whose inline-shape name needs
# to change from "image1.jpg" to "image2.jpg", and two PearsonLogo
# images (one in the "default" footer, one in the "first page" footer)
# whose inline-shape name needs to change from "image2.png" to
# "image1.png". The embedded pictures / relationships themselves are
# untouched - this only renames the shapes.

function Set-InlineShapeName($ishp, $newName) {
    try {
        $ishp.Name = $newName
    } catch {
        # Re-fetch through the shape's own Range and retry, as a fresh
        # handle sometimes needs to be re-addressed after the story was
        # first touched.
        $ishp.Range.InlineShapes.Item(1).Name = $newName
    }
}

$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

# wdHeaderFooterPrimary = 1, wdHeaderFooterFirstPage = 2
$headerKinds = 1, 2
$footerKinds = 1, 2

foreach ($kind in $headerKinds) {
    $hdr = $sec.Headers.Item($kind)
    if ($hdr.Exists) {
        $count = $hdr.Range.InlineShapes.Count
        for ($i = 1; $i -le $count; $i++) {
            $ishp = $hdr.Range.InlineShapes.Item($i)
            if ($ishp.AlternativeText -eq "BTec_Logo-Orange") {
                Set-InlineShapeName $ishp "image2.jpg"
            }
        }
    }
}

foreach ($kind in $footerKinds) {
    $ftr = $sec.Footers.Item($kind)
    if ($ftr.Exists) {
        $count = $ftr.Range.InlineShapes.Count
        for ($i = 1; $i -le $count; $i++) {
            $ishp = $ftr.Range.InlineShapes.Item($i)
            if ($ishp.AlternativeText -like "*PearsonLogo.png") {
                Set-InlineShapeName $ishp "image1.png"
            }
        }
    }
}
